$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column R (general_college_subjects.arts),
# shifting the existing R:AE block to U:AH.
$ws.Range("R1:T1").EntireColumn.Insert()

# New column headers (row 1)
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New column data (row 2)
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 0

# Normalize existing descriptive text values to lowercase
$ws.Range("D2").Value = "not considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "very important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
